# Replace the Formpack skeleton template body with the new,
# fully-unformatted localisation-token template.
#
# Word's Range.InsertXML() REPLACES the target range's contents
# wholesale with the supplied WordprocessingML, so every old
# paragraph mark / run / pPr / rPr under the skeleton paragraphs
# is discarded along with the old text, and the freshly inserted
# paragraphs carry no formatting overrides at all -- matching the
# plain <w:p><w:r><w:t>...</w:t></w:r></w:p> shape in the target.

$d = $word.ActiveDocument

# One new (unformatted) paragraph per array entry, in document order.
$newParagraphs = @(
    '<w:p><w:r><w:t>{{INS t.__PACK_ID__.title}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{INS t.__PACK_ID__.description}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{INS t.__PACK_ID__.section.person.title}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{INS t.__PACK_ID__.person.name.label}}: {{INS person.name}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{INS t.__PACK_ID__.person.birthDate.label}}: {{INS person.birthDate}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{INS t.__PACK_ID__.person.email.label}}: {{INS person.email}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{INS t.__PACK_ID__.person.website.label}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{IF person.website}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{LINK ({ url: person.website, label: person.website })}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{END-IF}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{INS t.__PACK_ID__.section.contacts.title}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{FOR c IN contacts}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{ALIAS contactName INS $c.name}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{ALIAS contactPhone INS $c.phone}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{INS t.__PACK_ID__.contacts.entry.label}} {{INS $idx}}: {{*contactName}} ({{*contactPhone}})</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{FOR detail IN [$c.phone, $c.relation].filter(Boolean)}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>- {{INS $detail}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{END-FOR detail}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{END-FOR c}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{INS t.__PACK_ID__.section.medications.title}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{FOR m IN medications}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{INS $idx}}. {{INS $m.name}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{IF $m.dosage}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{INS t.__PACK_ID__.medications.dosage.label}}: {{INS $m.dosage}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{END-IF}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{IF $m.schedule}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{INS t.__PACK_ID__.medications.schedule.label}}: {{INS $m.schedule}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{END-IF}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{END-FOR m}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{INS t.__PACK_ID__.section.diagnoses.title}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{IF diagnoses.formatted}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{INS diagnoses.formatted}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{END-IF}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{INS t.__PACK_ID__.section.symptoms.title}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{INS symptoms}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{INS t.__PACK_ID__.section.allergies.title}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{INS allergies}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{INS t.__PACK_ID__.section.doctor.title}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{INS t.__PACK_ID__.doctor.name.label}}: {{INS doctor.name}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{INS t.__PACK_ID__.doctor.phone.label}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{IF doctor.phone}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{LINK ({ url: ''tel:'' + doctor.phone, label: doctor.phone })}}</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>{{END-IF}}</w:t></w:r></w:p>'
)

$newBody = [string]::Join('', $newParagraphs)

$newBodyXml = '<?xml version="1.0" standalone="yes"?>' + 
    '<?mso-application progid="Word.Document"?>' + 
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + 
    '<pkg:part pkg:name="/word/document.xml" ' + 
        'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + 
    '<pkg:xmlData>' + 
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + 
    '<w:body>' + $newBody + '</w:body>' + 
    '</w:document>' + 
    '</pkg:xmlData></pkg:part></pkg:package>'

# InsertXML on the whole document Content range swaps in the new body
# (the trailing sectPr, which Content does not include, is left alone).
$d.Content.InsertXML($newBodyXml)

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"
